# Replace the sample appointment row (row 2) with new data.
# Columns: A=Stomatolog, B=Nume, C=Prenume, D=Ora, E=Nr Telefon,
#          F=Procedura, G=Data, H=Ora Data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 ("03847474") and G2 ("5/24/2024") look like a number / a date, so force
# them to be stored as literal text before assigning — otherwise Excel would
# silently reinterpret them (dropping the leading zero / turning the date
# into a serial number) instead of keeping the exact string.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"

$ws.Range("A2").Value = "Iulia"
$ws.Range("B2").Value = "Iulia"
$ws.Range("C2").Value = "Iulia"
$ws.Range("D2").Value = "12:00"
$ws.Range("E2").Value = "03847474"
$ws.Range("F2").Value = "Înălbire"
$ws.Range("G2").Value = "5/24/2024"
$ws.Range("H2").Value = "12:005/24/2024"
